# maj template comment à la fin
# Move the "Comment" column (and its related rows 2-5 metadata) from
# column J to the end (column P), shifting SamplePortion..ResultUnit
# left by one column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: column headers
$ws.Range("J1").Value = "SamplePortion"
$ws.Range("K1").Value = "SamplePortionUnit"
$ws.Range("L1").Value = "CookingTime"
$ws.Range("M1").Value = "Temperature"
$ws.Range("N1").Value = "Result"
$ws.Range("O1").Value = "ResultUnit"
$ws.Range("P1").Value = "Comment"

# Row 2: French labels
$ws.Range("J2").Value = "# Prise d'essai"
$ws.Range("K2").Value = "# Unité de mesure de la prise d’essai"
$ws.Range("L2").Value = "#TempsCuisson"
$ws.Range("M2").Value = "#Temperature"
$ws.Range("N2").Value = "# Résultat"
$ws.Range("O2").Value = "# Unité du résultat"
$ws.Range("P2").Value = "# Commentaire"

# Row 3: type hints
$ws.Range("J3").Value = "#float"
$ws.Range("K3").Value = "#string"
$ws.Range("L3").Value = "#integer"
$ws.Range("M3").Value = "#integer,`n  unit:celsius"
$ws.Range("N3").Value = "#float"
$ws.Range("O3").Value = "#string"
$ws.Range("P3").Value = "#string"

# Row 4: format hints
$ws.Range("J4").Value = "# format: nombre décimal, ne pas spécifier d'unité"
$ws.Range("K4").Value = "# format: texte"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = "# format: nombre décimal ou NA"
$ws.Range("O4").Value = "# format: texte"
$ws.Range("P4").Value = "# format: texte libre"

# Row 5: example values
$ws.Range("J5").Value = "# ex: 2.5"
$ws.Range("K5").Value = "# ex: mg"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = "# 409.935 ou NA"
$ws.Range("O5").Value = "# ex: mg/ml"
$ws.Range("P5").Value = ""
